# Updates cryptos list values (price + 1h volume change) to the latest
# scraped snapshot, and fixes the swapped Algorand / FraxShare rows (42-43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.669.50'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.846.98'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.035'
$ws.Range('E4').Value = '  +0.72%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.22'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.030'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4386'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3785'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07390'
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8810'
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.49'
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').Value = '1.874.91'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.486'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.700'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07120'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '84.77'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009060'
$ws.Range('E18').Value = '  +0.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.031'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.44'
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').Value = '27.718.83'
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.284'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.26'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('D24').Value = '2.089.28'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.051'
$ws.Range('E25').Value = '  +6.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.64'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.69'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.988'
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.326'
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.60'
$ws.Range('E30').Value = '  +1.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09065'
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7708'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.205'
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.997'
$ws.Range('E34').Value = '  +4.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.550'
$ws.Range('E35').Value = '  +0.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.032'
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.149'
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01972'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05252'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.849'
$ws.Range('E40').Value = '  +2.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5177'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.869'
$ws.Range('E42').Value = '  +2.96%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1665'
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.688'
$ws.Range('E44').Value = '  +1.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.95'
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.74'
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.032'
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06563'
$ws.Range('E48').Value = '  +2.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.699'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4688'
$ws.Range('E50').Value = '  +0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.878'
$ws.Range('E51').Value = '  -0.85%  '
